$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D/E on this sheet hold free-form text (prices with local
# formatting, padded percentages) that must stay text even when the
# digits happen to parse as a number, so force text format on D before
# writing to avoid Excel silently re-typing e.g. "65.80" -> 65.8.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.574.94'
$ws.Range('E2').Value = '  -2.56%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.656.04'
$ws.Range('E3').Value = '  -4.31%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.94'
$ws.Range('E5').Value = '  -2.48%  '
$ws.Range('E6').Value = '  -2.29%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.04'
$ws.Range('E8').Value = '  -0.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.261'
$ws.Range('E9').Value = '  -2.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0619'
$ws.Range('E10').Value = '  -2.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0881'
$ws.Range('E11').Value = '  -1.69%  '
$ws.Range('E12').Value = '  -4.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.642.96'
$ws.Range('E13').Value = '  -5.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.14'
$ws.Range('E14').Value = '  -2.79%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.565'
$ws.Range('E15').Value = '  +0.21%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.80'
$ws.Range('E16').Value = '  -2.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.580.46'
$ws.Range('E17').Value = '  -2.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '240.94'
$ws.Range('E18').Value = '  -2.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0728'
$ws.Range('E19').Value = '  -3.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.54'
$ws.Range('E20').Value = '  -4.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.46'
$ws.Range('E22').Value = '  -3.92%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.37'
$ws.Range('E23').Value = '  -3.31%  '
$ws.Range('E24').Value = '  -2.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.85'
$ws.Range('E25').Value = '  -2.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.20'
$ws.Range('E26').Value = '  -3.99%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.27'
$ws.Range('E27').Value = '  -2.85%  '
$ws.Range('E28').Value = '  +0.13%  '
$ws.Range('E29').Value = '  -2.43%  '
$ws.Range('E30').Value = '  +0.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0502'
$ws.Range('E31').Value = '  -2.66%  '
$ws.Range('E32').Value = '  -2.76%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.452.47'
$ws.Range('E33').Value = '  -2.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.11'
$ws.Range('E34').Value = '  -4.95%  '
$ws.Range('E35').Value = '  -4.81%  '
$ws.Range('E36').Value = '  -1.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.920'
$ws.Range('E37').Value = '  -5.84%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.572'
$ws.Range('E38').Value = '  -5.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0171'
$ws.Range('E39').Value = '  -2.81%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '69.64'
$ws.Range('E40').Value = '  -0.58%  '
$ws.Range('E41').Value = '  -3.62%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('E43').Value = '  -3.92%  '
$ws.Range('E44').Value = '  -3.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.794'
$ws.Range('E45').Value = '  -0.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.800.64'
$ws.Range('E46').Value = '  -4.14%  '
$ws.Range('E47').Value = '  -1.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.52'
$ws.Range('E48').Value = '  -2.14%  '
$ws.Range('E49').Value = '  -5.96%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.102'
$ws.Range('E50').Value = '  -1.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.82'
$ws.Range('E51').Value = '  -4.28%  '
